# p3.xlsx — widen columns B and E, and refresh the computed sample values
# (rows 2-9, columns A-E). Row 1 (the header-ish constants 10/20/40/60/100)
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15.08984375
$ws.Columns.Item(5).ColumnWidth = 14.453125

# --- Row 2 -----------------------------------------------------------------
$ws.Range("A2").Value = -0.56202855404733643
$ws.Range("B2").Value = -0.74939940164032248
$ws.Range("C2").Value = 0.3554238887755945
$ws.Range("D2").Value = -0.35236708129572647
$ws.Range("E2").Value = -0.050884515603718185

# --- Row 3 -----------------------------------------------------------------
$ws.Range("A3").Value = -0.63593667521452135
$ws.Range("B3").Value = -0.82048077239253547
$ws.Range("C3").Value = 0.37721155866329692
$ws.Range("D3").Value = -0.37603217112587239
$ws.Range("E3").Value = -0.042825639724322608

# --- Row 4 -----------------------------------------------------------------
$ws.Range("A4").Value = 0.021976819779258729
$ws.Range("B4").Value = 0.028613310289766779
$ws.Range("C4").Value = -0.05247976060291637
$ws.Range("D4").Value = 0.039139611834479275
$ws.Range("E4").Value = 0.0060304726236641067

# --- Row 5 -----------------------------------------------------------------
$ws.Range("A5").Value = 0.022112109671254036
$ws.Range("B5").Value = 0.029223298562533653
$ws.Range("C5").Value = -0.053252880674375745
$ws.Range("D5").Value = 0.039716114344862427
$ws.Range("E5").Value = 0.0062609664618908474

# --- Row 6 -----------------------------------------------------------------
$ws.Range("A6").Value = 0.93050538302384844
$ws.Range("B6").Value = -0.396570190400583
$ws.Range("C6").Value = 0.74726395435573811
$ws.Range("D6").Value = -0.78659797009336652
$ws.Range("E6").Value = -0.99388905282225759

# --- Row 7 -----------------------------------------------------------------
$ws.Range("A7").Value = 0.9923654818025508
$ws.Range("B7").Value = -0.46360157096328186
$ws.Range("C7").Value = 0.83259608810925734
$ws.Range("D7").Value = -0.87331615342773194
$ws.Range("E7").Value = -1.089258851688542

# --- Row 8 -----------------------------------------------------------------
$ws.Range("A8").Value = -0.090925075188737423
$ws.Range("B8").Value = -0.0025860958626791803
$ws.Range("C8").Value = -0.02394903184583097
$ws.Range("D8").Value = 0.052888343205285793
$ws.Range("E8").Value = 0.066367335699086788

# --- Row 9 -----------------------------------------------------------------
$ws.Range("A9").Value = -0.092149235455031744
$ws.Range("B9").Value = -0.0024740593617054636
$ws.Range("C9").Value = -0.024721875441250592
$ws.Range("D9").Value = 0.053648071154433784
$ws.Range("E9").Value = 0.067615058342343354
